# Update "想去人数" (column F) values on the "展览" and "全部类型" sheets
# to reflect newer counts, as scraped at commit 456a3b4.

$wb = $excel.ActiveWorkbook

# Mapping of row number -> new value for column F
$updates = @{
    7  = 2495
    9  = 126
    11 = 1263
    15 = 1127
    16 = 314
    18 = 24
    19 = 25
    21 = 67
    22 = 83
    23 = 216
    25 = 254
}

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($row in $updates.Keys) {
        $ws.Range("F$row").Value = $updates[$row]
    }
}
